$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New component stats table (rows 2-52), replacing the old rows 2-54.
$data = @(
    @{Row=2; Name="Menu"; Bugs=17; Features=5; Total=22},
    @{Row=3; Name="DataGrid"; Bugs=5; Features=6; Total=11},
    @{Row=4; Name="Tree"; Bugs=10; Features=0; Total=10},
    @{Row=5; Name="Nav"; Bugs=8; Features=1; Total=9},
    @{Row=6; Name="Combobox"; Bugs=4; Features=4; Total=8},
    @{Row=7; Name="Popover"; Bugs=6; Features=2; Total=8},
    @{Row=8; Name="Dialog"; Bugs=6; Features=1; Total=7},
    @{Row=9; Name="Tooltip"; Bugs=4; Features=1; Total=5},
    @{Row=10; Name="Table"; Bugs=4; Features=1; Total=5},
    @{Row=11; Name="TagPicker"; Bugs=4; Features=1; Total=5},
    @{Row=12; Name="Dropdown"; Bugs=3; Features=2; Total=5},
    @{Row=13; Name="Toolbar"; Bugs=3; Features=1; Total=4},
    @{Row=14; Name="Virtualizer"; Bugs=4; Features=0; Total=4},
    @{Row=15; Name="TeachingPopover"; Bugs=3; Features=0; Total=3},
    @{Row=16; Name="Skeleton"; Bugs=2; Features=1; Total=3},
    @{Row=17; Name="MessageBar"; Bugs=3; Features=0; Total=3},
    @{Row=18; Name="Calendar Compat"; Bugs=3; Features=0; Total=3},
    @{Row=19; Name="Slider"; Bugs=0; Features=2; Total=2},
    @{Row=20; Name="Portal"; Bugs=2; Features=0; Total=2},
    @{Row=21; Name="Toast"; Bugs=1; Features=1; Total=2},
    @{Row=22; Name="Drawer"; Bugs=2; Features=0; Total=2},
    @{Row=23; Name="Accordion"; Bugs=2; Features=0; Total=2},
    @{Row=24; Name="FluentProvider"; Bugs=0; Features=2; Total=2},
    @{Row=25; Name="Tabs"; Bugs=2; Features=0; Total=2},
    @{Row=26; Name="DatePicker"; Bugs=0; Features=2; Total=2},
    @{Row=27; Name="Switch"; Bugs=2; Features=0; Total=2},
    @{Row=28; Name="Avatar"; Bugs=2; Features=0; Total=2},
    @{Row=29; Name="FocusTrapZone"; Bugs=1; Features=0; Total=1},
    @{Row=30; Name="Popup"; Bugs=1; Features=0; Total=1},
    @{Row=31; Name="Button"; Bugs=0; Features=1; Total=1},
    @{Row=32; Name="Tag"; Bugs=1; Features=0; Total=1},
    @{Row=33; Name="DatePickerCompat"; Bugs=0; Features=1; Total=1},
    @{Row=34; Name="Label"; Bugs=1; Features=0; Total=1},
    @{Row=35; Name="Image"; Bugs=0; Features=1; Total=1},
    @{Row=36; Name="Input"; Bugs=0; Features=1; Total=1},
    @{Row=37; Name="AvatarGroup"; Bugs=1; Features=0; Total=1},
    @{Row=38; Name="SearchBox"; Bugs=1; Features=0; Total=1},
    @{Row=39; Name="Carousel"; Bugs=1; Features=0; Total=1},
    @{Row=40; Name="List"; Bugs=1; Features=0; Total=1},
    @{Row=41; Name="InfoLabel"; Bugs=1; Features=0; Total=1},
    @{Row=42; Name="Persona"; Bugs=1; Features=0; Total=1},
    @{Row=43; Name="SplitButton"; Bugs=1; Features=0; Total=1},
    @{Row=44; Name="MenuItem"; Bugs=1; Features=0; Total=1},
    @{Row=45; Name="Badge"; Bugs=1; Features=0; Total=1},
    @{Row=46; Name="Checkbox"; Bugs=1; Features=0; Total=1},
    @{Row=47; Name="Rating"; Bugs=0; Features=0; Total=0},
    @{Row=48; Name="ColorPicker"; Bugs=0; Features=0; Total=0},
    @{Row=49; Name="Pickers"; Bugs=0; Features=0; Total=0},
    @{Row=50; Name="Keytip"; Bugs=0; Features=0; Total=0},
    @{Row=51; Name="Segment"; Bugs=0; Features=0; Total=0},
    @{Row=52; Name="SpinButton"; Bugs=0; Features=0; Total=0}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.Name
    $ws.Cells.Item($r, 2).Value = $item.Bugs
    $ws.Cells.Item($r, 3).Value = $item.Features
    $ws.Cells.Item($r, 4).Value = $item.Total
}

# Remove the two now-obsolete trailing rows (previously ColorPicker, FloatingLabelInput)
$ws.Range("A53:D54").Delete()
